$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combine the old A2 (name), A3 (type) and A4 (text) into a single python-tuple-like string in A2
$ws.Range("A2").Value = "('Flooded Strand', ['Land', '{T}, Pay 1 life, Sacrifice Flooded Strand: Search your library for a Plains or Island card, put it onto the battlefield, then shuffle your library.'])"

# Remove the now-unused rows 3 and 4 entirely so the sheet dimension shrinks back to A1:A2
$ws.Range("A3:A4").EntireRow.Delete()
